$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("PIB por Estado")

# Grab the chart object so we can reposition it after the row shuffle.
$co = $ws.ChartObjects().Item(1)

# Remove the footnote / source rows at the bottom of the sheet
# ("Cifras Preliminares:", "/p1 ...", "Cifras Revisadas:", "/r1 ...",
#  "Fuentes:", "/f1 ...") which lived in rows 14-19.
$ws.Rows("14:19").Delete() | Out-Null

# Remove the title rows at the top of the sheet ("INEGI..." source line,
# "Banco de Informacion..." line, and the blank spacer row) which lived
# in rows 1-3. This shifts the "Entidad Federativa"/"Periodo 2020" header
# (old row 4) up to row 1, and the data rows (old 5-11) up to rows 2-8.
$ws.Rows("1:3").Delete() | Out-Null

# The chart anchor needs to move up to follow the header row back to the
# top of the sheet (it used to start at row 4 / sit over rows 4-18, now it
# should start at row 1 / sit over rows 1-11). The chart's column-based
# horizontal placement doesn't change, only its vertical placement.
$co.Top = 0
$co.Height = 127.9

# Leave the selection where the author ended up after the edits.
$ws.Range("A21").Select() | Out-Null

Write-Host "Done"
